# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits the leading word + trailing space of each title's first run
# ("First " / "Third ") into two separate runs ("First"/"Third" and a
# lone " ") instead of leaving them consolidated into a single run.
#
# Re-assigning a Characters() sub-range's Text to itself is enough to
# force the writer to materialize it as its own run (with a plain,
# untouched <a:rPr/>) without altering any character formatting.

$p = $ppt.ActivePresentation

function Split-LeadingWordRun($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    $title = $slide.Shapes.Item(1)
    $tr = $title.TextFrame.TextRange
    $full = $tr.Text
    $spaceIdx = $full.IndexOf(" ")
    if ($spaceIdx -lt 0) { return }

    # 1-based, length-1 character range covering just the word before
    # the space (e.g. "First" out of "First slide").
    $word = $tr.Characters(1, $spaceIdx)
    $word.Text = $word.Text

    # The single space character right after the word.
    $space = $tr.Characters($spaceIdx + 1, 1)
    $space.Text = $space.Text
}

Split-LeadingWordRun 1
Split-LeadingWordRun 3
